# Bill of materials update:
#  - add a new "Raspberry PI 4B case" component row (row 4)
#  - retarget the Raspberry PI 4B designator to "Single Board Computer"
#  - add category/version-style designator names ("Wireless transceivers", etc.)
#  - keep all existing rows + hyperlinks intact, just shifted down by one

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remember every existing hyperlink (address + optional literal display text) ---
# before-state refs, keyed by their current (pre-insert) cell address
$oldLinks = @(
    @{ Ref = "G7";  Address = "https://www.digikey.com/en/products/detail/american-opto-plus-led/L314ED/13677700?s=N4IgTCBcDaIDIGYCMAWAogERAXQL5A"; Display = "https://www.digikey.com/en/products/detail/american-opto-plus-led/L314ED/13677700?s=N4IgTCBcDaIDIGYCMAWAogERAXQL5A" },
    @{ Ref = "G10"; Address = "https://www.digikey.com/en/products/detail/kemet/c1206s104k5racauto/10232834"; Display = "https://www.digikey.com/en/products/detail/kemet/c1206s104k5racauto/10232834" },
    @{ Ref = "G5";  Address = "https://oshpark.com/"; Display = $null },
    @{ Ref = "G9";  Address = "https://www.digikey.com/en/products/detail/panasonic-electronic-components/ECE-A1EKS100B/2689095"; Display = $null },
    @{ Ref = "G6";  Address = "https://www.digikey.com/en/products/detail/raspberry-pi/SC0339L/12339165?gad_source=1&gad_campaignid=20243136172&gbraid=0AAAAADrbLliTb25Erv6NiD1YUmheOuryx&gclid=Cj0KCQjwvJHIBhCgARIsAEQnWlDQ6JN4eg4O2KvihoihYaPILt0_g-A7s0a8kHtxZmLVeqAIM-Sh0BoaAvo5EALw_wcB&gclsrc=aw.ds"; Display = "https://www.digikey.com/en/products/detail/raspberry-pi/SC0339L/12339165?gad_source=1&gad_campaignid=20243136172&gbraid=0AAAAADrbLliTb25Erv6NiD1YUmheOuryx&gclid=Cj0KCQjwvJHIBhCgARIsAEQnWlDQ6JN4eg4O2KvihoihYaPILt0_g-A7s0a8kHtxZmLVeqAIM-Sh0BoaAvo5EALw_wcB&gclsrc=aw.ds" },
    @{ Ref = "H7";  Address = "https://www.lvelectronics.com/details/item?itemid=AOP%20L314ED"; Display = $null },
    @{ Ref = "H9";  Address = "https://www.mouser.com/ProductDetail/Panasonic/ECE-A1EKS100B?qs=rMMd5vBiahpNSBku%2FPA9KQ%3D%3D&srsltid=AfmBOoqQ4U271ZrVJxuPMHVWnjSZpLGUJXgJG3eSAxUWLyhrYK65hiV8"; Display = $null },
    @{ Ref = "H10"; Address = "https://www.mouser.com/ProductDetail/KEMET/C1206S104K5RACAUTO?qs=ds50AKTGxA9C%252BOGFjvBCPg%3D%3D&srsltid=AfmBOoqQltyUeWiBl72NEX0atiJAYoEHcEyVBj5XktYi3QvV_cm_2CP3"; Display = $null },
    @{ Ref = "H6";  Address = "https://www.staples.com/verbatim-premium-64gb-microsdxc-memory-card-with-adapter-class-10-uhs-i-v10-44084/product_2478898"; Display = $null },
    @{ Ref = "I6";  Address = "https://www.cdw.com/product/sandisk-extreme-flash-memory-card-64-gb-microsdxc-uhs-i/7316165?pfm=srh"; Display = $null },
    @{ Ref = "I3";  Address = "https://www.mouser.com/ProductDetail/Raspberry-Pi/SC01949?qs=T%252BzbugeAwjjISb%252BwlagpRw%3D%3D"; Display = $null },
    @{ Ref = "G4";  Address = "https://protosupplies.com/product/nrf24l01palna-2-4ghz-rf-wireless-module/"; Display = $null },
    @{ Ref = "H4";  Address = "https://www.elecbee.com/en/product-detail/1100-meter-long-distance-nrf24l01-pa-lna-wireless-module-with-antenna-module_73275?utm_term=&utm_campaign=&utm_source=adwords&utm_medium=ppc&hsa_acc=9958698819&hsa_cam=23146566611&hsa_grp=187297116859&hsa_ad=779498650524&hsa_src=g&hsa_tgt=pla-2511885727437&hsa_kw=&hsa_mt=&hsa_net=adwords&hsa_ver=3&gad_source=1&gad_campaignid=23146566611&gbraid=0AAAAADGHwHYfKF0iMEZtoqyx1rWjzHw1t&gclid=Cj0KCQiA5abIBhCaARIsAM3-zFXbXR-v0fRxqHBuXN-ZMrGgtzU7CUVsQzAEYIulzVkVFA95q0AEql8aAkk6EALw_wcB"; Display = "https://www.elecbee.com/en/product-detail/1100-meter-long-distance-nrf24l01-pa-lna-wireless-module-with-antenna-module_73275?utm_term=&utm_campaign=&utm_source=adwords&utm_medium=ppc&hsa_acc=9958698819&hsa_cam=23146566611&hsa_grp=187297116859&hsa_ad=779498650524&hsa_src=g&hsa_tgt=pla-2511885727437&hsa_kw=&hsa_mt=&hsa_net=adwords&hsa_ver=3&gad_source=1&gad_campaignid=23146566611&gbraid=0AAAAADGHwHYfKF0iMEZtoqyx1rWjzHw1t&gclid=Cj0KCQiA5abIBhCaARIsAM3-zFXbXR-v0fRxqHBuXN-ZMrGgtzU7CUVsQzAEYIulzVkVFA95q0AEql8aAkk6EALw_wcB" },
    @{ Ref = "I4";  Address = "https://openelab.io/products/nrf24l01-transceiver-transceiver-module?srsltid=AfmBOorW9lfupEOXFS3kU72MZQgoLedxFpkylsIVPOrEJnmP53qSnV0F"; Display = $null },
    @{ Ref = "G8";  Address = "https://www.digikey.com/en/products/detail/stackpole-electronics-inc/CF14JT270R/1741362"; Display = $null },
    @{ Ref = "H8";  Address = "https://www.mouser.com/ProductDetail/YAGEO/CFR25SJT-26-270R?qs=sGAEpiMZZMsPqMdJzcrNwjRWpSA1Ui8vfWaBDCu2IB3jHFqFn4CeHg%3D%3D"; Display = $null }
)

# drop every hyperlink up front -- the row insert below does not renumber them for us
$ws.Cells.Hyperlinks.Delete()

# --- Insert the new component row above the old row 4 (NRF24 row), shifting rows 4-26 down to 5-27 ---
$ws.Rows("4:4").Insert()
$ws.Rows("4:4").RowHeight = 15

# --- Re-point the Raspberry Pi designator (row 3) at the new "Single Board Computer" category ---
$ws.Range("A3").Value = "Single Board Computer"

# --- Populate the new row 4: Raspberry PI 4B case ---
$ws.Range("A4").Value = "case for Single Board Computer"
$ws.Range("B4").Value = "Raspberry PI 4B case"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "pc"
$ws.Range("E4").Value = 16
$ws.Range("F4").Formula = "=C4*E4"
$ws.Range("G4").Value = "https://www.digikey.com/en/products/detail/edatec/ED-PI5CASE-BS/21769634"
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = ""

# --- Re-create every pre-existing hyperlink at its shifted location (rows >= 4 move down by one) ---
foreach ($link in $oldLinks) {
    $ref = $link.Ref
    $col = $ref -replace '[0-9]', ''
    $row = [int]($ref -replace '[A-Za-z]', '')
    if ($row -ge 4) { $row = $row + 1 }
    $newRef = "$col$row"
    if ($link.Display) {
        $ws.Hyperlinks.Add($ws.Range($newRef), $link.Address, "", "", $link.Display)
    } else {
        $ws.Hyperlinks.Add($ws.Range($newRef), $link.Address)
    }
}

# --- Add the two brand-new hyperlinks for the Single Board Computer rows ---
$ws.Hyperlinks.Add($ws.Range("G3"), "https://www.digikey.com/en/products/detail/raspberry-pi/SC1642/24627138")
$ws.Hyperlinks.Add($ws.Range("G4"), "https://www.digikey.com/en/products/detail/edatec/ED-PI5CASE-BS/21769634")

# --- Match the author's final selection state ---
$ws.Range("F17").Select()

$wb.Save()
